$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Fix the timestamp already stored in row 190 (A190) and fill in its
# previously-blank precio_num / fecha_dia cells (E190 / F190).
# ---------------------------------------------------------------------------
$ws.Range("A190").Value = 45986.43551783565
$ws.Range("E190").Value = 32.91
$ws.Range("F190").Value = 45986

# ---------------------------------------------------------------------------
# Append the new EVOWHEY PROTEIN price-log rows 191-198 (fully populated),
# then row 199 whose precio_num / fecha_dia are intentionally left blank,
# matching the scraper's latest (incomplete) run.
# ---------------------------------------------------------------------------
$newRows = @(
    @{ Row = 191; A = 45986.48644931713; E = 32.91; F = 45986 },
    @{ Row = 192; A = 45986.49722497685; E = 32.91; F = 45986 },
    @{ Row = 193; A = 45987.42156263889; E = 32.91; F = 45987 },
    @{ Row = 194; A = 45987.42373547453; E = 32.91; F = 45987 },
    @{ Row = 195; A = 45987.42436534722; E = 32.91; F = 45987 },
    @{ Row = 196; A = 45987.42464135417; E = 32.91; F = 45987 },
    @{ Row = 197; A = 45987.42620861111; E = 32.91; F = 45987 },
    @{ Row = 198; A = 45987.42665641203; E = 32.91; F = 45987 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = "EVOWHEY PROTEIN"
    $ws.Cells.Item($row, 3).Value = "2Kg"
    $ws.Cells.Item($row, 4).Value = "32,91€"
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
}

# Row 199: same product columns, but precio_num / fecha_dia stay blank.
$ws.Cells.Item(199, 1).Value = 45987.42761194793
$ws.Cells.Item(199, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(199, 3).Value = "2Kg"
$ws.Cells.Item(199, 4).Value = "32,91€"

# ---------------------------------------------------------------------------
# Number formats: column A uses the workbook's date-time format, column F
# (when populated) uses the plain date format - same styles already used
# throughout the rest of the sheet.
# ---------------------------------------------------------------------------
$ws.Range("A190:A199").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("F190:F198").NumberFormat = "YYYY-MM-DD"
